# Add 7 new rows (156-162) of requisition data to the worksheet, replicating
# the formatting/styles of the last existing row (155) and filling in the
# new values (matching the upstream commit "Add files via upload").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 155
$newRowsCount = 7

# Duplicate the formatting (styles/number formats) of the last row into each
# of the new rows first, so number formats / styles match the rest of the
# table exactly (copy formats only, not values).
for ($i = 1; $i -le $newRowsCount; $i++) {
    $targetRow = $lastRow + $i
    $ws.Range("A$lastRow`:P$lastRow").Copy()
    $ws.Range("A$targetRow`:P$targetRow").PasteSpecial(-4122)
    # Match the row height used throughout the rest of the table
    # (ht="12.75" customHeight="1").
    $ws.Rows($targetRow).RowHeight = 12.75
}

# Row data: EMPRD, EMPRD_DESC, EMPRD_UF, REQ_CDG, REQ_DATA, OF_CDG, OF_DATA,
# INSUMO_CDG, INSUMO_DESC, INSUMO_STATUS, INSUMO_UNIDADE, QTD_PED,
# ITEM_PRCUNTPED, PRCTTL_INSUMO, FORNECEDOR_CDG, FORNECEDOR_DESC
#
# FORNECEDOR_CDG ("O") values are zero-padded numeric-looking codes that
# must remain text. Rather than typing them in (which would make Excel
# coerce them to numbers and drop the leading zeros), we point at an
# existing row elsewhere in the sheet that already holds that exact
# supplier code as text, and copy that single cell (value + style) down -
# OSrcRow identifies a donor row with the matching FORNECEDOR_CDG/DESC text.
$rows = @(
    @{ Row=156; A=2317; B="LUIZ ALBERTO HESS BORGES"; C="SP"; D=63; E=46013.4672055093; F=81382; G=46013; H="E.02.0065"; I="PAPELÃO COM LOGO OSBORNE"; J="Apto"; K="UN"; L=2;   M=90;   N=180;   OSrcRow=3;   P="CASA PEDROSO2648864-" },
    @{ Row=157; A=2317; B="LUIZ ALBERTO HESS BORGES"; C="SP"; D=63; E=46013.4672055093; F=81381; G=46013; H="E.03.0150"; I="BOTA  DE SEGURANÇA MSA  NOBUCK MARLUVAS CADARÇO"; J="Apto"; K="PAR"; L=1; M=190;  N=190;   OSrcRow=97;  P="GALPÃO DAS FERRAMENT" },
    @{ Row=158; A=2317; B="LUIZ ALBERTO HESS BORGES"; C="SP"; D=63; E=46013.4672055093; F=81381; G=46013; H="E.04.0892"; I="DISCO DE CORTE AÇO INOX 4 1/2''"; J="Apto"; K="UN"; L=10;  M=4.15; N=41.5;  OSrcRow=97;  P="GALPÃO DAS FERRAMENT" },
    @{ Row=159; A=2317; B="LUIZ ALBERTO HESS BORGES"; C="SP"; D=63; E=46013.4672055093; F=81382; G=46013; H="H.11.0033"; I="AÇO CA50 8,0 MM - VARA"; J="Apto"; K="UN"; L=8; M=38.9; N=311.2; OSrcRow=3;   P="CASA PEDROSO2648864-" },
    @{ Row=160; A=2317; B="LUIZ ALBERTO HESS BORGES"; C="SP"; D=63; E=46013.4672055093; F=81382; G=46013; H="J.03.0015"; I="AREIA  - SACO GRANDE 20KG"; J="Apto"; K="SC"; L=400; M=5.8; N=2320; OSrcRow=3;   P="CASA PEDROSO2648864-" },
    @{ Row=161; A=2317; B="LUIZ ALBERTO HESS BORGES"; C="SP"; D=63; E=46013.4672055093; F=81382; G=46013; H="M.09.0023"; I="ARAME RECOZIDO FIO DUPLO TRANÇADO NØ 18"; J="Apto"; K="KG"; L=5; M=14.9; N=74.5; OSrcRow=3;   P="CASA PEDROSO2648864-" },
    @{ Row=162; A=2317; B="LUIZ ALBERTO HESS BORGES"; C="SP"; D=63; E=46013.4672055093; F=81383; G=46013; H="P3.02.0171"; I="RALO HEMISFÉRICO ( ABACAXI  ) - 100MM -"; J="Apto"; K="UN"; L=15; M=28; N=420; OSrcRow=140; P="ATLANTA" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.A
    $ws.Cells.Item($row, 2).Value2 = $r.B
    $ws.Cells.Item($row, 3).Value2 = $r.C
    $ws.Cells.Item($row, 4).Value2 = $r.D
    $ws.Cells.Item($row, 5).Value2 = $r.E
    $ws.Cells.Item($row, 6).Value2 = $r.F
    $ws.Cells.Item($row, 7).Value2 = $r.G
    $ws.Cells.Item($row, 8).Value2 = $r.H
    $ws.Cells.Item($row, 9).Value2 = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N

    # Copy the donor cell's value+style verbatim so the zero-padded code
    # is preserved as text without introducing any new style entry.
    $ws.Range("O$($r.OSrcRow)").Copy()
    $ws.Range("O$row").PasteSpecial(-4104)

    $ws.Cells.Item($row, 16).Value2 = $r.P
}
